$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cast")

# Force the new row to be stored as text (matches the existing rows, which
# are all plain text values even when they look numeric/date-like), so
# Excel doesn't reinterpret "2024-10-02" as a date serial or "7" as a number.
$ws.Range("A11:E11").NumberFormat = "@"

$ws.Range("A11").Value = "2024-10-02"
$ws.Range("B11").Value = "AZD"
$ws.Range("C11").Value = "7"
$ws.Range("D11").Value = "asd"
$ws.Range("E11").Value = "rqwe"
